$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: swap in the "invalid"/"valid" status values (previously in row 6 area)
$ws.Range("A1").Value = "invalid"
$ws.Range("B1").Value = "valid"

# Row 6: now holds what used to be in row 1 (the credential pair)
$ws.Range("A6").Value = "mngr251101"
$ws.Range("B6").Value = "dYrYhun"

# Update the selected cell in the sheet view
$ws.Range("F6").Select()
